# Swap the contents of columns B:AB between specific pairs of rows.
# Column A (the running "id" index) stays untouched; every other field
# (match id, teams, odds, etc.) for the two rows trades places.
#
# Value2 is used instead of Value because it reliably round-trips both
# numeric and text cell contents through this COM-interop runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(17, 18),
    @(73, 74),
    @(78, 79),
    @(173, 174)
)

# Columns B (2) through AB (28)
$firstCol = 2
$lastCol = 28

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
